# Generate Report for Handback
# Update the handoff/handback timestamps for the
# "7621df27-4b0e-4f4d-b275-68936b0ec56b" file row across the
# Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the
#     7621df27... row (row 3) moves from 10:51:49 to 10:52:43.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-30 10:52:43"

# --- zh-cn sheet: Correspond Handoff Datetime / Correspond Handback
#     DateTime for the 7621df27... row (row 3).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-30 10:52:39"
$wsZhCn.Range("K3").Value = "2016-08-30 10:52:57"

# --- de-de sheet: Correspond Handoff Datetime / Correspond Handback
#     DateTime for the 7621df27... row (row 3). The Handoff Datetime
#     value mirrors the Overview sheet's "Latest HO Xliff Generate
#     Date" for this row, so it is set to the same new value here.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-30 10:52:43"
$wsDeDe.Range("K3").Value = "2016-08-30 10:53:12"
